$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K ("Дата передачи показаний" / Reading submission date),
# shifting the existing date columns (Дата поверки ... Дата контрольных показаний)
# one column to the right (K:P -> L:Q).
$ws.Columns("K").Insert()

# New column header
$ws.Range("K1").Value = "Дата передачи показаний"

# The other date-like columns in this sheet store their values as plain text
# (e.g. "2021-01-20"), not as real dates - match that by forcing Text format
# before writing the values so Excel doesn't auto-convert them to date serials.
$ws.Range("K2:K10").NumberFormat = "@"
$ws.Range("K2:K6").Value = "2021-12-20"
$ws.Range("K7:K10").Value = "2021-12-21"

# Restore the new column's width (close to the original authored width).
$ws.Columns("K").ColumnWidth = 20
